# Apply cryptos list update (price/volume refresh) as produced by the
# scheduled GitHub Actions job. Mirrors the OOXML diff cell-by-cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" cells contain plain numeric-looking text (e.g. "1.002",
# "0.00001339", "1.0000") that Excel would otherwise silently reinterpret
# as a real number (dropping trailing zeros / using scientific notation).
# Force those specific cells to Text format first so the literal string
# is preserved exactly, just like the original inline-string cells.
$numericLookingPriceCells = @(
    "D4","D5","D6","D7","D8","D9","D10","D11",
    "D12","D13","D14","D15","D16","D18","D19","D20",
    "D21","D22","D23","D25","D26","D27","D29","D30",
    "D31","D32","D34","D35","D36","D37","D38","D39",
    "D40","D41","D42","D43","D44","D45","D46","D47",
    "D49","D50"
)
foreach ($addr in $numericLookingPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "24.377.28"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "1.666.71"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "312.89"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "0.3942"
$ws.Range("E7").Value = "  +1.68%  "
$ws.Range("D8").Value = "0.3922"
$ws.Range("E8").Value = "  +2.22%  "
$ws.Range("D9").Value = "52.26"
$ws.Range("E9").Value = "  +6.26%  "
$ws.Range("D10").Value = "1.395"
$ws.Range("E10").Value = "  +4.02%  "
$ws.Range("D11").Value = "0.9999"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "0.08575"
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("D13").Value = "24.44"
$ws.Range("E13").Value = "  +3.42%  "
$ws.Range("D14").Value = "7.311"
$ws.Range("E14").Value = "  +3.28%  "
$ws.Range("D15").Value = "7.966"
$ws.Range("E15").Value = "  +7.15%  "
$ws.Range("D16").Value = "0.00001339"
$ws.Range("E16").Value = "  +5.15%  "
$ws.Range("D17").Value = "1.660.54"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").Value = "94.89"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "0.06989"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").Value = "20.57"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "6.994"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("D22").Value = "1.0000"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").Value = "13.75"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("D24").Value = "24.387.57"
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "3.059"
$ws.Range("E25").Value = "  +14.71%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "2.424"
$ws.Range("E26").Value = "  +3.86%  "
$ws.Range("D27").Value = "22.53"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "143.09"
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").Value = "5.433"
$ws.Range("E30").Value = "  +3.59%  "
$ws.Range("D31").Value = "8.052"
$ws.Range("E31").Value = "  -7.72%  "
$ws.Range("D32").Value = "2.538"
$ws.Range("E32").Value = "  +3.42%  "
$ws.Range("D33").Value = "1.846.20"
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").Value = "1.059"
$ws.Range("E34").Value = "  +10.87%  "
$ws.Range("D35").Value = "0.08242"
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("D36").Value = "0.03041"
$ws.Range("E36").Value = "  +4.77%  "
$ws.Range("D37").Value = "6.946"
$ws.Range("E37").Value = "  -2.22%  "
$ws.Range("D38").Value = "11.13"
$ws.Range("E38").Value = "  +12.31%  "
$ws.Range("D39").Value = "0.2765"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("D40").Value = "0.09237"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("D41").Value = "0.7702"
$ws.Range("E41").Value = "  +2.03%  "
$ws.Range("D42").Value = "13.79"
$ws.Range("E42").Value = "  +6.43%  "
$ws.Range("D43").Value = "1.454"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").Value = "16.54"
$ws.Range("E44").Value = "  +4.56%  "
$ws.Range("D45").Value = "0.7125"
$ws.Range("E45").Value = "  +4.06%  "
$ws.Range("D46").Value = "2.539"
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("D47").Value = "4.139"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").Value = "0.08424"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("D50").Value = "136.58"
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("E51").Value = "  +1.55%  "
